$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D. This shifts the existing D:K (FY2017..FY2010
# + trailing blank) data right by one column to E:L, matching the diff's column shift.
$ws.Columns("D:D").Insert()

# The freshly inserted column D cells come back with the default/general style.
# Copy number formats from column E (which now holds what used to be column D,
# i.e. the correct per-row style) back onto the new column D so every cell keeps
# the same look (date style for the header row, #,##0 style for the data rows).
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)

# Populate the new column D with the newest fiscal year (FY2018, period ending
# 2018-12-31 = serial 43465) financial figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 1556300
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = -21800
$ws.Range("D17").Value = 350300
$ws.Range("D18").Value = 1206000
$ws.Range("D20").Value = -629700
$ws.Range("D21").Value = 634400
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 576300
$ws.Range("D24").Value = 117400
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 458900
$ws.Range("D27").Value = 444600
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 9200
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 629700
$ws.Range("D33").Value = 453800
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 453800
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 665700
$ws.Range("D42").Value = 707500
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 17200
$ws.Range("D48").Value = 456300
$ws.Range("D49").Value = 2865700
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 89400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 47877300
$ws.Range("D57").Value = "NA"
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 906800
$ws.Range("D62").Value = "NA"
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 41343400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 244100
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 1284800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 6289800
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 453800
$ws.Range("D83").Value = 58100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 589300
$ws.Range("D91").Value = -87900
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -400
$ws.Range("D96").Value = -257900
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -539500
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 49400
